# Generate Report for Archive
#
# 1. "Ready for handoff" -> "In Translation" everywhere it is used
#    (Overview!E2:E4 + F2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. The localization-status columns that held that text get narrower
#    (autofit-style shrink now that the text is shorter):
#       Overview columns E & F, zh-cn column C, de-de column C
#    old stored width ~17.216  ->  new stored width ~13.410

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:E4").Value = $newStatus
$wsOverview.Range("F2:F4").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus

# Narrow the affected columns to match the new (shorter) content width.
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
